﻿# edit.ps1
# Reproduces:
#  1. Slide 16's table switches to table style {70634070-4277-4C46-92F8-6A94669162AC}
#     (was {8C5D474E-1256-452B-9732-A22C46AED82E}).
#  2. The deck's two theme parts swap identity: the theme backing the slide
#     master (Integral colours) ends up carrying the Office default palette,
#     while the theme backing the notes master (Office default palette) ends
#     up carrying the Integral palette.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 16 -------------------------------------------
$slide = $p.Slides.Item(16)
$tableShape = $slide.Shapes.Item(3)
$table = $tableShape.Table
$table.ApplyStyle("{70634070-4277-4C46-92F8-6A94669162AC}")

# --- 2. Swap the two themes' colour palettes -------------------------------
# theme1.xml backs the slide master and currently holds the "Integral" palette.
# theme2.xml backs the notes master and currently holds the "Office" palette.
# Swap their 12 theme colours (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
# so each part ends up with the other's colours.

$slideMasterColors = $p.SlideMaster.Theme.ThemeColorScheme
$notesMasterColors = $p.NotesMaster.Theme.ThemeColorScheme

$integralRgb = @(0, 16777215, 5332805, 13754083, 3722137, 3646819, 2412774, 38860, 13611854, 10915127, 2465643, 158642)
$officeRgb   = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

for ($i = 1; $i -le 12; $i++) {
    $slideMasterColors.Item($i).RGB = $officeRgb[$i - 1]
    $notesMasterColors.Item($i).RGB = $integralRgb[$i - 1]
}
